$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 30 (pushing the previous
# row 30 and everything below it down by one row, dimension A1:R79 -> A1:R80).
$ws.Rows("30").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 44533
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = 100112032
$ws.Range("G30").Value = "Zapallo italiano"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 250
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 6500
$ws.Range("M30").Value = 6300
$ws.Range("N30").Value = "$/caja 60 unidades"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 105
$ws.Range("Q30").Value = 60
$ws.Range("R30").Value = "Hortaliza"
